$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row appended at the bottom of the sheet (row 86)
$row = 86

# A86: date value, formatted like the other date cells in column A (style index 1) -
# copy the style from the row above so we reuse the existing date style instead of
# minting a new one, then overwrite with the new value.
$ws.Range("A85").Copy($ws.Range("A$row"))
$ws.Range("A$row").Value = 45447.2916666667

# B-F: plain numeric values (volume, high, low, open, close)
$ws.Range("B$row").Value = 6100
$ws.Range("C$row").Value = 0.720000028610229
$ws.Range("D$row").Value = 0.704999983310699
$ws.Range("E$row").Value = 0.704999983310699
$ws.Range("F$row").Value = 0.720000028610229

# G86: adj_close, stored as text (matches existing text-typed column in sheet)
$ws.Range("G$row").Value = "'0.720000028610229"
$ws.Range("G$row").ClearFormats()

# H86: ticker, stored as text
$ws.Range("H$row").Value = "BWZ.MI"
